$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.838.29"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "3.054.59"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'548.15"
$ws.Range("E5").Value = "  +2.08%  "

$ws.Range("D6").Value = "'135.98"
$ws.Range("E6").Value = "  -0.60%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.052.46"
$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  +1.06%  "

$ws.Range("D10").Value = "'6.20"
$ws.Range("E10").Value = "  +0.06%  "

$ws.Range("D11").Value = "'0.149"
$ws.Range("E11").Value = "  -3.93%  "

$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").Value = "'35.28"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("D14").Value = "'0.0000224"
$ws.Range("E14").Value = "  +1.18%  "

$ws.Range("D15").Value = "3.543.51"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").Value = "62.881.31"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.060.93"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.110"
$ws.Range("E18").Value = "  -2.40%  "

$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("D20").Value = "'486.28"
$ws.Range("E20").Value = "  +3.77%  "

$ws.Range("D21").Value = "'13.42"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").Value = "'7.11"
$ws.Range("E23").Value = "  +1.97%  "

$ws.Range("D24").Value = "'81.77"
$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("D25").Value = "'12.26"
$ws.Range("E25").Value = "  +1.54%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("E27").Value = "  +1.67%  "

$ws.Range("D28").Value = "'7.92"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").Value = "'1.95"
$ws.Range("E30").Value = "  +4.54%  "

$ws.Range("D31").Value = "'26.02"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").Value = "'5.75"
$ws.Range("E33").Value = "  +5.62%  "

$ws.Range("D34").Value = "'2.38"
$ws.Range("E34").Value = "  +3.73%  "

$ws.Range("D35").Value = "'55.56"

$ws.Range("D36").Value = "'5.94"
$ws.Range("E36").Value = "  +0.26%  "

$ws.Range("D37").Value = "'460.08"
$ws.Range("E37").Value = "  -3.98%  "

$ws.Range("D38").Value = "3.186.27"
$ws.Range("E38").Value = "  -1.99%  "

$ws.Range("D39").Value = "'0.0810"
$ws.Range("E39").Value = "  +2.57%  "

$ws.Range("D40").Value = "'0.0391"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("E41").Value = "  +1.94%  "

$ws.Range("D42").Value = "'8.20"
$ws.Range("E42").Value = "  +1.16%  "

$ws.Range("D43").Value = "'2.48"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").Value = "'26.52"
$ws.Range("E44").Value = "  +5.36%  "

$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("E46").Value = "  -0.41%  "

$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("D48").Value = "'2.01"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("D49").Value = "'116.73"
$ws.Range("E49").Value = "  -5.08%  "

$ws.Range("D50").Value = "0.0₃0498"
$ws.Range("E50").Value = "  -3.93%  "

$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "'1.28"
$ws.Range("E51").Value = "  +3.96%  "

$ws.Range("D2:E51").Style = "Normal"
